$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 10 (Objetivos:) B/C value - previously showed the wrong text
$ws.Range("B10").Value = 'Geral: Formar profissionais em nível superior com capacidade de conhecer a sequência dos procedimentos de análise química de interesse ambiental. Executar procedimentos de análises volumétricas bem como interpretar, avaliar e criticar os resultados obtidos. Específicos: Formar profissionais em nível superior com capacidade de conhecer as etapas da sequência analítica. Conhecer os métodos de preparação de amostras a serem analisadas. Compreender as bases teóricas da química analítica quantitativa de interesse ambiental. Executar procedimentos de análises químicas volumétricas (volumetria por neutralização, volumetria por precipitação, volumetria por oxi-redução e volumetria por complexação), análises ambientais de DBO, DQO, OD e turbidez, bem como interpretar, avaliar e criticar os resultados obtidos.'
$ws.Range("C10").Value = 'Geral: Formar profissionais em nível superior com capacidade de conhecer a sequência dos procedimentos de análise química de interesse ambiental. Executar procedimentos de análises volumétricas bem como interpretar, avaliar e criticar os resultados obtidos. Específicos: Formar profissionais em nível superior com capacidade de conhecer as etapas da sequência analítica. Conhecer os métodos de preparação de amostras a serem analisadas. Compreender as bases teóricas da química analítica quantitativa de interesse ambiental. Executar procedimentos de análises químicas volumétricas (volumetria por neutralização, volumetria por precipitação, volumetria por oxi-redução e volumetria por complexação), análises ambientais de DBO, DQO, OD e turbidez, bem como interpretar, avaliar e criticar os resultados obtidos.'

# Insert new row 13 to hold the "Docentes responsaveis:" value, shifting rows 13-23 down to 14-24
$ws.Rows("13:13").Insert()
$ws.Range("A13").Clear()
$ws.Range("B3").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B13").Value = '8855158 - Morun Bernardino Neto'
$ws.Range("C13").Value = '8855158 - Morun Bernardino Neto'

# Fix row 14 (Programa resumido:) B/C value - previously "Semestral"
$ws.Range("B14").Value = 'Fundamentos da análise titulométrica (preparação de amostras e cálculos). Análises titulométricas por precipitação, neutralização, complexação e oxirredução. Análises ambientais de DBO, DQO, OD e turbidez.'
$ws.Range("C14").Value = 'Fundamentos da análise titulométrica (preparação de amostras e cálculos). Análises titulométricas por precipitação, neutralização, complexação e oxirredução. Análises ambientais de DBO, DQO, OD e turbidez.'

# Fix row 16 (Programa:) B/C value - previously wrongly showed a date
$ws.Range("B16").Value = '- Fundamentos de análise titulométrica e cálculos em análise titulométrica.- Titulometria de neutralização: fundamentos, indicadores de titulação, curvas de titulação ácido base.- Titulometria complexométrica: fundamentos, complexometria com EDTA.- Titulação de oxiredução: fundamentos e principais indicadores- Titulometria de precipitação: fundamentos, indicadores, argentimetria.- Análises quantitativas em solo, água, alimentos etc.- Análises ambientais: DBO, DQO, OD e turbidez'
$ws.Range("C16").Value = '- Fundamentos de análise titulométrica e cálculos em análise titulométrica.- Titulometria de neutralização: fundamentos, indicadores de titulação, curvas de titulação ácido base.- Titulometria complexométrica: fundamentos, complexometria com EDTA.- Titulação de oxiredução: fundamentos e principais indicadores- Titulometria de precipitação: fundamentos, indicadores, argentimetria.- Análises quantitativas em solo, água, alimentos etc.- Análises ambientais: DBO, DQO, OD e turbidez'

# Fix row 19 (Metodo:) B/C value - previously showed wrong text
$ws.Range("B19").Value = 'O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais.'
$ws.Range("C19").Value = 'O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais.'

# Fix row 20 (Criterio:) B/C value
$ws.Range("B20").Value = 'Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.'
$ws.Range("C20").Value = 'Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.'

# Fix row 21 (Norma de recuperacao:) B/C value
$ws.Range("B21").Value = 'Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 e estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.'
$ws.Range("C21").Value = 'Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 e estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.'

# Fix row 22 (Bibliografia:) B/C value - previously showed the recovery-norm text; now shows the bibliography text
$ws.Range("B22").Value = 'Harris, D.C. EXPLORANDO A QUÍMICA ANALÍTICA, 4ª edição, LTC, Rio de Janeiro – RJ, 2011Skoog, D.A., Holler, F.J. e Nieman, T.A., PRINCÍPIOS DE ANÁLISE INSTRUMENTAL, 5ª ed., Bookman, Porto Alegre, 2002.Mendham,J., Denney, R.C., Barnes, J.D. e Thomas, M., Vogel: ANÁLISE QUÍMICA QUANTITATIVA, 6ª ed., Livros Técnicos e Científicos, Rio de Janeiro -RJ, 2002.'
$ws.Range("C22").Value = 'Harris, D.C. EXPLORANDO A QUÍMICA ANALÍTICA, 4ª edição, LTC, Rio de Janeiro – RJ, 2011Skoog, D.A., Holler, F.J. e Nieman, T.A., PRINCÍPIOS DE ANÁLISE INSTRUMENTAL, 5ª ed., Bookman, Porto Alegre, 2002.Mendham,J., Denney, R.C., Barnes, J.D. e Thomas, M., Vogel: ANÁLISE QUÍMICA QUANTITATIVA, 6ª ed., Livros Técnicos e Científicos, Rio de Janeiro -RJ, 2002.'

# Column width cleanup: column A only spans col 1 (was erroneously min=1 max=2)
$ws.Columns("A").ColumnWidth = 30.7109375
$ws.Columns("B").ColumnWidth = 60.7109375
$ws.Columns("C").ColumnWidth = 60.7109375
